# Updates cached price/profit figures across the Leve-profit worksheets
# (columns H-N: currentAveragePrice[NQ/HQ], LevePrice[NQ/HQ], LeveProfit[NQ/HQ])
# to the latest scheduled-runner snapshot values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 95.2
$ws.Range("I6").Value = 95.2
$ws.Range("K6").Value = 285.6
$ws.Range("M6").Value = -173.6

# Row 12
$ws.Range("H12").Value = 318.33334
$ws.Range("I12").Value = 318.33334
$ws.Range("K12").Value = 318.33334
$ws.Range("M12").Value = -148.33334

# Row 17
$ws.Range("H17").Value = 390.39
$ws.Range("J17").Value = 390.39
$ws.Range("L17").Value = 1171.17
$ws.Range("N17").Value = -1507.17

# Row 51
$ws.Range("H51").Value = 20364.705
$ws.Range("I51").Value = 9837.5
$ws.Range("J51").Value = 29722.223
$ws.Range("K51").Value = 9837.5
$ws.Range("L51").Value = 29722.223
$ws.Range("M51").Value = -9353.5
$ws.Range("N51").Value = -30690.223

# Row 103
$ws.Range("H103").Value = 1754.5714
$ws.Range("J103").Value = 1862.8334
$ws.Range("L103").Value = 5588.5002
$ws.Range("N103").Value = -6760.5002

# Row 112
$ws.Range("H112").Value = 3632.4644
$ws.Range("J112").Value = 3632.4644
$ws.Range("L112").Value = 10897.3932
$ws.Range("N112").Value = -13113.3932

# Row 137
$ws.Range("H137").Value = 7582147
$ws.Range("I137").Value = 13893398
$ws.Range("K137").Value = 41680194
$ws.Range("M137").Value = -41677644

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 9286.471
$ws.Range("I61").Value = 4861.273
$ws.Range("J61").Value = 17399.334
$ws.Range("K61").Value = 4861.273
$ws.Range("L61").Value = 17399.334
$ws.Range("M61").Value = -4649.273
$ws.Range("N61").Value = -17823.334

# Row 122
$ws.Range("H122").Value = 3263.5
$ws.Range("I122").Value = 2559.6667
$ws.Range("K122").Value = 7679.000100000001
$ws.Range("M122").Value = -5229.000100000001

# Row 130
$ws.Range("H130").Value = 148910.56
$ws.Range("J130").Value = 148910.56
$ws.Range("L130").Value = 148910.56
$ws.Range("N130").Value = -158950.56

# Row 136
$ws.Range("H136").Value = 9286.471
$ws.Range("I136").Value = 4861.273
$ws.Range("J136").Value = 17399.334
$ws.Range("K136").Value = 14583.819
$ws.Range("L136").Value = 52198.00199999999
$ws.Range("M136").Value = -12033.819
$ws.Range("N136").Value = -57298.00199999999

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 2170.8
$ws.Range("I58").Value = 2170.8
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 2170.8
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -1967.8
$ws.Range("N58").ClearContents()

# Row 86
$ws.Range("H86").Value = 7914.778
$ws.Range("J86").Value = 6866.8
$ws.Range("L86").Value = 6866.8
$ws.Range("N86").Value = -9112.799999999999

# Row 89
$ws.Range("H89").Value = 7914.778
$ws.Range("J89").Value = 6866.8
$ws.Range("L89").Value = 34334
$ws.Range("N89").Value = -45566

# Row 132
$ws.Range("H132").Value = 137816.05
$ws.Range("I132").Value = 95429.10000000001
$ws.Range("K132").Value = 286287.3
$ws.Range("M132").Value = -283757.3

# Row 136
$ws.Range("H136").Value = 2170.8
$ws.Range("I136").Value = 2170.8
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 6512.400000000001
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -3962.400000000001
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 2858.08
$ws.Range("J2").Value = 5032.5713
$ws.Range("L2").Value = 30195.4278
$ws.Range("N2").Value = -30421.4278

# Row 31
$ws.Range("H31").Value = 200
$ws.Range("I31").Value = 200
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 600
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -312
$ws.Range("N31").ClearContents()

# Row 34
$ws.Range("H34").Value = 2142.6428
$ws.Range("I34").Value = 79
$ws.Range("J34").Value = 3289.111
$ws.Range("K34").Value = 237
$ws.Range("L34").Value = 9867.332999999999
$ws.Range("M34").Value = -153
$ws.Range("N34").Value = -10035.333

# Row 39
$ws.Range("H39").Value = 1895.4166
$ws.Range("J39").Value = 4400
$ws.Range("L39").Value = 13200
$ws.Range("N39").Value = -13788

# Row 55
$ws.Range("H55").Value = 924.1667
$ws.Range("I55").Value = 770.8
$ws.Range("J55").Value = 1033.7142
$ws.Range("K55").Value = 2312.4
$ws.Range("L55").Value = 3101.1426
$ws.Range("M55").Value = -2135.4
$ws.Range("N55").Value = -3455.1426

# Row 121
$ws.Range("H121").Value = 27778384
$ws.Range("I121").Value = 541.8889
$ws.Range("J121").Value = 111111910
$ws.Range("K121").Value = 1625.6667
$ws.Range("L121").Value = 333335730
$ws.Range("M121").Value = -315.6667000000002
$ws.Range("N121").Value = -333338350

# Row 131
$ws.Range("H131").Value = 16674565
$ws.Range("I131").Value = 83334740
$ws.Range("J131").Value = 9522.625
$ws.Range("K131").Value = 250004220
$ws.Range("L131").Value = 28567.875
$ws.Range("M131").Value = -249999180
$ws.Range("N131").Value = -38647.875

$ws = $wb.Worksheets.Item("GSM")
# Row 3
$ws.Range("H3").Value = 1965.2142
$ws.Range("I3").Value = 1551.0834
$ws.Range("J3").Value = 4450
$ws.Range("K3").Value = 1551.0834
$ws.Range("L3").Value = 4450
$ws.Range("M3").Value = -1435.0834
$ws.Range("N3").Value = -4682

# Row 69
$ws.Range("H69").Value = 17100
$ws.Range("I69").Value = 15000
$ws.Range("J69").Value = 18500
$ws.Range("K69").Value = 15000
$ws.Range("L69").Value = 18500
$ws.Range("M69").Value = -14251
$ws.Range("N69").Value = -19998

# Row 70
$ws.Range("H70").Value = 16641.715
$ws.Range("I70").Value = 14956.25
$ws.Range("K70").Value = 14956.25
$ws.Range("M70").Value = -14686.25

# Row 72
$ws.Range("H72").Value = 17100
$ws.Range("I72").Value = 15000
$ws.Range("J72").Value = 18500
$ws.Range("K72").Value = 45000
$ws.Range("L72").Value = 55500
$ws.Range("M72").Value = -41256
$ws.Range("N72").Value = -62988

# Row 73
$ws.Range("H73").Value = 16641.715
$ws.Range("I73").Value = 14956.25
$ws.Range("K73").Value = 14956.25
$ws.Range("M73").Value = -14020.25

# Row 80
$ws.Range("H80").Value = 3841.0588
$ws.Range("I80").Value = 3556.4167
$ws.Range("J80").Value = 4524.2
$ws.Range("K80").Value = 3556.4167
$ws.Range("L80").Value = 4524.2
$ws.Range("M80").Value = -2558.4167
$ws.Range("N80").Value = -6520.2

# Row 83
$ws.Range("H83").Value = 3841.0588
$ws.Range("I83").Value = 3556.4167
$ws.Range("J83").Value = 4524.2
$ws.Range("K83").Value = 17782.0835
$ws.Range("L83").Value = 22621
$ws.Range("M83").Value = -12790.0835
$ws.Range("N83").Value = -32605

# Row 102
$ws.Range("H102").Value = 1944.3125
$ws.Range("I102").Value = 1968.9231
$ws.Range("K102").Value = 1968.9231
$ws.Range("M102").Value = -346.9231

# Row 123
$ws.Range("H123").Value = 24000
$ws.Range("J123").Value = 24000
$ws.Range("L123").Value = 24000
$ws.Range("N123").Value = -28900

# Row 132
$ws.Range("H132").Value = 3459.3704
$ws.Range("I132").Value = 1856.95
$ws.Range("K132").Value = 5570.85
$ws.Range("M132").Value = -3040.85

$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 605.5454999999999
$ws.Range("I55").Value = 841.4
$ws.Range("J55").Value = 409
$ws.Range("K55").Value = 841.4
$ws.Range("L55").Value = 409
$ws.Range("M55").Value = -668.4
$ws.Range("N55").Value = -755

# Row 136
$ws.Range("H136").Value = 4046.4849
$ws.Range("I136").Value = 2703.0356
$ws.Range("K136").Value = 8109.1068
$ws.Range("M136").Value = -5559.1068

$ws = $wb.Worksheets.Item("WVR")
# Row 61
$ws.Range("H61").Value = 34825.168
$ws.Range("I61").Value = 32262.75
$ws.Range("K61").Value = 32262.75
$ws.Range("M61").Value = -31970.75

# Row 122
$ws.Range("H122").Value = 3162.25
$ws.Range("I122").Value = 2924.6667
$ws.Range("K122").Value = 8774.000100000001
$ws.Range("M122").Value = -6324.000100000001

# Row 126
$ws.Range("H126").Value = 3069.3635
$ws.Range("I126").Value = 1529.6111
$ws.Range("J126").Value = 9998.25
$ws.Range("K126").Value = 4588.8333
$ws.Range("L126").Value = 29994.75
$ws.Range("M126").Value = -2118.8333
$ws.Range("N126").Value = -34934.75

# Row 135
$ws.Range("H135").Value = 113141.336
$ws.Range("J135").Value = 113141.336
$ws.Range("L135").Value = 113141.336
$ws.Range("N135").Value = -123281.336
